# "fixed export and fixing maps"
# - Rename the worksheet from "1" to "სენაკი"
# - Remove the "(მოსახლეობის აღწერის შედეგებით)" caption row (row 2)
# - Collapse the 1989/2002/2014 columns down to just the 2014 column
#   by deleting the two now-unwanted year columns (B:C), leaving the
#   former "2014" column (D) as the new column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the extra census-year columns (1989, 2002); the 2014 column shifts left.
$ws.Range("B:C").Delete()

# Drop the now-unused subtitle row "(მოსახლეობის აღწერის შედეგებით)".
$ws.Rows.Item(2).Delete()

# Rename sheet to match the municipality name.
$ws.Name = "სენაკი"
